# Implementing multiple files support
#
# The "single attachment" columns (AK..AS: Totale imposta, Riferimento
# normativo, Condizioni di pagamento, Modalita, Importo, Istituto
# finanziario, Codice IBAN, Nome dell'allegato, Formato) are removed,
# since an invoice can now come with multiple attached files, each
# represented as its own data row (columns A..AJ) instead of extra
# trailing columns tacked on to the first row.
#
# Two further invoice rows (sheet rows 4 and 5, in addition to the
# existing row 3... no: rows 3-5 are new) are appended below the existing
# data row 2, reusing the same column layout (A..AJ) as the header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove the obsolete AK:AS columns (header row 1 + data row 2) ---
# Shifting left collapses/removes these trailing cells entirely (there is
# nothing further right to take their place), which also shrinks the
# sheet's used range/dimension back down to column AJ.
$ws.Range("AK1:AS2").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)

# --- Step 2: append new data rows 3-5 with additional invoice records ---
#
# Every column in these new rows is stored as text (even the numeric-
# looking ones, matching how the existing data row 2 stores things like
# " 1400.00" as text) except column A, which is a real number. Force
# text typing by temporarily applying a "@" (text) number format to the
# destination cells before writing the values - this stops Excel from
# silently re-parsing strings like " 11" or " 2500.00" into numbers -
# then clear the formatting back off afterwards so the cells end up with
# no explicit style (matching the source rows).

$textRange = $ws.Range("B3:AJ5")
$textRange.NumberFormat = "@"

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = " IT01879020517"
$ws.Range("C3").Value = " 11"
$ws.Range("D3").Value = " FPR12"
$ws.Range("E3").Value = " 0000000"
$ws.Range("F3").Value = " monservicesrls@pec.it"
$ws.Range("G3").Value = " 1"
$ws.Range("H3").Value = " FSNMRC74C14F109F"
$ws.Range("I3").Value = " Fattura a saldo per lavori edili svolti c/o vs cantiere sito in Carpignano"
$ws.Range("J3").Value = " RF01 (ordinario)"
$ws.Range("K3").Value = " LE"
$ws.Range("L3").Value = " 73040"
$ws.Range("M3").Value = " IT"
$ws.Range("N3").Value = " TD01 (fattura)"
$ws.Range("O3").Value = " EUR"
$ws.Range("P3").Value = " 2500.00"
$ws.Range("Q3").Value = " IT04879980755"
$ws.Range("R3").Value = " 04879980755"
$ws.Range("S3").Value = " MON SERVICE SOCIETA' A RESPONSABILITA' LIMITATA SEMPLIFICATA"
$ws.Range("T3").Value = " TENUTA SPECCHIA, 2"
$ws.Range("U3").Value = " 1.00"
$ws.Range("V3").Value = " NR"
$ws.Range("W3").Value = " 2500.00"
$ws.Range("X3").Value = " 2500.00"
$ws.Range("Y3").Value = " 0.00"
$ws.Range("Z3").Value = " N6.3 (inversione contabile - subappalto nel settore edile)"
$ws.Range("AA3").Value = " 0.00"
$ws.Range("AB3").Value = " N6.3 (inversione contabile - subappalto nel settore edile)"
$ws.Range("AC3").Value = " 2500.00"
$ws.Range("AD3").Value = " 0.00"
$ws.Range("AE3").Value = " Reverse charge subappalto nel settore edile Art. 17, c. 6 lett. a), DPR 633/72"
$ws.Range("AF3").Value = " TP02 (pagamento completo)"
$ws.Range("AG3").Value = " MP05 (bonifico)"
$ws.Range("AH3").Value = " 2022-04-11 (11 Aprile 2022)"
$ws.Range("AI3").Value = " 2500.00"
$ws.Range("AJ3").Value = " INTESA SAN PAOLO"

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = " IT01879020517"
$ws.Range("C4").Value = " 12"
$ws.Range("D4").Value = " FPR12"
$ws.Range("E4").Value = " 0000000"
$ws.Range("F4").Value = " monservicesrls@pec.it"
$ws.Range("G4").Value = " 1"
$ws.Range("H4").Value = " FSNMRC74C14F109F"
$ws.Range("I4").Value = " Fattura a saldo per lavori edili di ristrutturazione effettuati per vs conto c/o"
$ws.Range("J4").Value = " RF01 (ordinario)"
$ws.Range("K4").Value = " LE"
$ws.Range("L4").Value = " 73040"
$ws.Range("M4").Value = " IT"
$ws.Range("N4").Value = " TD01 (fattura)"
$ws.Range("O4").Value = " EUR"
$ws.Range("P4").Value = " 6000.00"
$ws.Range("Q4").Value = " IT04879980755"
$ws.Range("R4").Value = " 04879980755"
$ws.Range("S4").Value = " MON SERVICE SOCIETA' A RESPONSABILITA' LIMITATA SEMPLIFICATA"
$ws.Range("T4").Value = " TENUTA SPECCHIA, 2"
$ws.Range("U4").Value = " 1.00"
$ws.Range("V4").Value = " NR"
$ws.Range("W4").Value = " 6000.00"
$ws.Range("X4").Value = " 6000.00"
$ws.Range("Y4").Value = " 0.00"
$ws.Range("Z4").Value = " N6.3 (inversione contabile - subappalto nel settore edile)"
$ws.Range("AA4").Value = " 0.00"
$ws.Range("AB4").Value = " N6.3 (inversione contabile - subappalto nel settore edile)"
$ws.Range("AC4").Value = " 6000.00"
$ws.Range("AD4").Value = " 0.00"
$ws.Range("AE4").Value = " Reverse charge subappalto nel settore edile Art. 17, c. 6 lett. a), DPR 633/72"
$ws.Range("AF4").Value = " TP02 (pagamento completo)"
$ws.Range("AG4").Value = " MP05 (bonifico)"
$ws.Range("AH4").Value = " 2022-04-11 (11 Aprile 2022)"
$ws.Range("AI4").Value = " 6000.00"
$ws.Range("AJ4").Value = " INTESA SAN PAOLO"

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = " IT01879020517"
$ws.Range("C5").Value = " 15"
$ws.Range("D5").Value = " FPR12"
$ws.Range("E5").Value = " 0000000"
$ws.Range("F5").Value = " monservicesrls@pec.it"
$ws.Range("G5").Value = " 1"
$ws.Range("H5").Value = " FSNMRC74C14F109F"
$ws.Range("I5").Value = " Fattura di cappotto, rasato, smontaggio e montaggio marmi c/o immobile"
$ws.Range("J5").Value = " RF01 (ordinario)"
$ws.Range("K5").Value = " LE"
$ws.Range("L5").Value = " 73040"
$ws.Range("M5").Value = " IT"
$ws.Range("N5").Value = " TD01 (fattura)"
$ws.Range("O5").Value = " EUR"
$ws.Range("P5").Value = " 8100.00"
$ws.Range("Q5").Value = " IT04879980755"
$ws.Range("R5").Value = " 04879980755"
$ws.Range("S5").Value = " MON SERVICE SOCIETA' A RESPONSABILITA' LIMITATA SEMPLIFICATA"
$ws.Range("T5").Value = " TENUTA SPECCHIA, 2"
$ws.Range("U5").Value = " 1.00"
$ws.Range("V5").Value = " NR"
$ws.Range("W5").Value = " 8100.00"
$ws.Range("X5").Value = " 8100.00"
$ws.Range("Y5").Value = " 0.00"
$ws.Range("Z5").Value = " N6.3 (inversione contabile - subappalto nel settore edile)"
$ws.Range("AA5").Value = " 0.00"
$ws.Range("AB5").Value = " N6.3 (inversione contabile - subappalto nel settore edile)"
$ws.Range("AC5").Value = " 8100.00"
$ws.Range("AD5").Value = " 0.00"
$ws.Range("AE5").Value = " Reverse charge subappalto nel settore edile Art. 17, c. 6 lett. a), DPR 633/72"
$ws.Range("AF5").Value = " TP02 (pagamento completo)"
$ws.Range("AG5").Value = " MP05 (bonifico)"
$ws.Range("AH5").Value = " 2022-04-29 (29 Aprile 2022)"
$ws.Range("AI5").Value = " 8100.00"
$ws.Range("AJ5").Value = " INTESA SAN PAOLO"

# Drop the temporary "@" text format again so these cells end up with no
# explicit style, same as the source data.
$textRange.ClearFormats()

# The original sheet styles the "line number" column (A) with the bold/
# centered header style (index 1) on every data row; match that for the
# new rows too.
$ws.Range("A3").Style = $ws.Range("A2").Style
$ws.Range("A4").Style = $ws.Range("A2").Style
$ws.Range("A5").Style = $ws.Range("A2").Style
